# Implements the "Penality Reward System" data trim:
#  - Sheet 1 "Weekly Quantity": remove 5 weeks of data (old rows 11-15),
#    shrinking the used range from A1:B38 down to A1:B33.
#  - Sheet 2 "Monthly Trend": update one value (B5: 20 -> 15) and remove
#    2 months of data (old rows 6-7), shrinking the used range from
#    A1:B14 down to A1:B12.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Weekly Quantity")
$ws1.Rows("11:15").Delete() | Out-Null

$ws2 = $wb.Worksheets.Item("Monthly Trend")
$ws2.Range("B5").Value = 15
$ws2.Rows("6:7").Delete() | Out-Null
